$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the player it documents.
$ws.Name = "Dhawal Kulkarni"

# Insert a new leading column for the match number; this pushes every
# existing column (teamName..result) one slot to the right (B..M) while
# keeping their values intact.
$ws.Columns.Item(1).Insert()

# Populate the new "matchNo" column.
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "27th"
